$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Update the refreshed query timestamps on the "data" sheet (column F) ---
$dataSheet.Range("F2").Value  = "2021-10-05 14:33:19.625345"
$dataSheet.Range("F3").Value  = "2021-10-05 14:33:19.625357"
$dataSheet.Range("F4").Value  = "2021-10-05 14:33:19.625360"
$dataSheet.Range("F5").Value  = "2021-10-05 14:33:19.625363"
$dataSheet.Range("F6").Value  = "2021-10-05 14:33:19.625366"
$dataSheet.Range("F7").Value  = "2021-10-05 14:33:19.625369"
$dataSheet.Range("F8").Value  = "2021-10-05 14:33:19.625371"
$dataSheet.Range("F9").Value  = "2021-10-05 14:33:19.625373"
$dataSheet.Range("F10").Value = "2021-10-05 14:33:19.625376"
$dataSheet.Range("F11").Value = "2021-10-05 14:33:19.625379"
$dataSheet.Range("F12").Value = "2021-10-05 14:33:19.625381"

# --- Add the new "metadata" sheet, placed after "data" ---
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = "metadata"

# Reuse the existing header style (bold + bordered) from the "data" sheet
# for the metadata header row and the leading id column.
$dataSheet.Range("B1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)

$dataSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Brugada syndrome"
$newSheet.Range("C2").Value = 60

# data_version must stay text ("0.34"), not be coerced to a number
$dVersion = $newSheet.Range("D2")
$dVersion.NumberFormat = "@"
$dVersion.Value = "0.34"
$dVersion.ClearFormats()

$newSheet.Range("E2").Value = "2020-06-01T04:28:49.182809Z"
$newSheet.Range("F2").Value = "2021-10-05 14:33:19.621300"
$newSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/60/?format=json"

# Keep "data" as the active sheet/selection, as in the original workbook
$dataSheet.Activate()
